{"js": "// Add \"\u041c\u0438\u0430\u0445 \u0422\u0430\u043a\u0431\u0438\u0440\" right after \"\u041c\u0430\u0437\u0443\u043c\u0434\u0435\u0440 \u0428\u043e\u0443\u0432\u0438\u043a\" in the student list\n// that lives inside the title-page text box. The text box content is\n// duplicated twice in the underlying OOXML (an mc:Choice DrawingML\n// rendition and an mc:Fallback VML rendition of the very same shape),\n// so the new paragraph has to be added to both copies to keep them in\n// sync - exactly what Word itself does when a user edits the text box\n// on screen.\n//\n// Office.js's Word.Paragraph/Shape object model in this host does not\n// expose the nested paragraphs that live inside a text box, so we\n// locate the host paragraph (the one whose run holds the drawing/\n// shape) and surgically patch its OOXML, then write it back with\n// Range.insertOoxml(..., Word.InsertLocation.replace). This only\n// touches the text box paragraph - nothing else in the document is\n// read or rewritten.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The run of text immediately preceding the insertion point, common to\n// both the mc:Choice and mc:Fallback copies of the text box.\nconst marker = \"<w:t>\u041c\u0430\u0437\u0443\u043c\u0434\u0435\u0440 \u0428\u043e\u0443\u0432\u0438\u043a</w:t></w:r></w:p>\";\n\n// The new paragraph being inserted, matching the bold / Russian-language\n// run formatting used by the rest of the student list.\nconst newParagraph =\n  \"<w:p>\" +\n  \"<w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val=\\\"ru-RU\\\" w:bidi=\\\"ar-IQ\\\"/></w:rPr></w:pPr>\" +\n  \"<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val=\\\"ru-RU\\\" w:bidi=\\\"ar-IQ\\\"/></w:rPr>\" +\n  \"<w:t>\u041c\u0438\u0430\u0445 \u0422\u0430\u043a\u0431\u0438\u0440</w:t>\" +\n  \"</w:r>\" +\n  \"</w:p>\";\n\nlet hostParagraph = null;\nlet hostOoxml = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const candidate = paragraphs.items[i];\n  const result = candidate.getOoxml();\n  await context.sync();\n  if (result.value.indexOf(marker) !== -1) {\n    hostParagraph = candidate;\n    hostOoxml = result.value;\n    break;\n  }\n}\n\nif (!hostParagraph) {\n  throw new Error('Could not find the paragraph containing \"\u041c\u0430\u0437\u0443\u043c\u0434\u0435\u0440 \u0428\u043e\u0443\u0432\u0438\u043a\".');\n}\n\nconst occurrences = hostOoxml.split(marker).length - 1;\nif (occurrences !== 2) {\n  throw new Error(\"Expected 2 occurrences (Choice + Fallback) of the marker, found \" + occurrences);\n}\n\nconst updatedOoxml = hostOoxml.split(marker).join(marker + newParagraph);\n\nhostParagraph.insertOoxml(updatedOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Add \"\u041c\u0438\u0430\u0445 \u0422\u0430\u043a\u0431\u0438\u0440\" right after \"\u041c\u0430\u0437\u0443\u043c\u0434\u0435\u0440 \u0428\u043e\u0443\u0432\u0438\u043a\" in the student list\n# that lives inside the title-page text box. The text box content is\n# duplicated twice in the underlying OOXML (an mc:Choice DrawingML\n# rendition and an mc:Fallback VML rendition of the very same shape),\n# so the new paragraph has to land in both copies to keep them in sync\n# - exactly what Word itself does when a user edits the text box on\n# screen.\n#\n# Shape.TextFrame.TextRange in this host resolves to the shape's host\n# paragraph in the main story rather than the nested text-box story, so\n# instead we locate that host paragraph by its OpenXML content and\n# surgically patch it with Range.InsertXML (the COM equivalent of\n# Office.js's Range.insertOoxml(..., Replace)), which replaces only\n# that paragraph's contents - nothing else in the document is touched.\n\n$d = $word.ActiveDocument\n\n$marker = \"<w:t>\u041c\u0430\u0437\u0443\u043c\u0434\u0435\u0440 \u0428\u043e\u0443\u0432\u0438\u043a</w:t></w:r></w:p>\"\n$newParagraph = \"<w:p>\" + `\n  \"<w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val=`\"ru-RU`\" w:bidi=`\"ar-IQ`\"/></w:rPr></w:pPr>\" + `\n  \"<w:r><w:rPr><w:b/><w:bCs/><w:lang w:val=`\"ru-RU`\" w:bidi=`\"ar-IQ`\"/></w:rPr>\" + `\n  \"<w:t>\u041c\u0438\u0430\u0445 \u0422\u0430\u043a\u0431\u0438\u0440</w:t>\" + `\n  \"</w:r>\" + `\n  \"</w:p>\"\n\n$hostRange = $null\n$hostXml = $null\n\n$paragraphs = $d.Paragraphs\nfor ($i = 1; $i -le $paragraphs.Count; $i++) {\n  $candidateRange = $paragraphs.Item($i).Range\n  $candidateXml = $candidateRange.WordOpenXML\n  if ($candidateXml.Contains($marker)) {\n    $hostRange = $candidateRange\n    $hostXml = $candidateXml\n    break\n  }\n}\n\nif ($hostRange -eq $null) {\n  throw 'Could not find the paragraph containing \"\u041c\u0430\u0437\u0443\u043c\u0434\u0435\u0440 \u0428\u043e\u0443\u0432\u0438\u043a\".'\n}\n\n$occurrences = ([regex]::Matches($hostXml, [regex]::Escape($marker))).Count\nif ($occurrences -ne 2) {\n  throw \"Expected 2 occurrences (Choice + Fallback) of the marker, found $occurrences\"\n}\n\n$updatedXml = $hostXml.Replace($marker, $marker + $newParagraph)\n$hostRange.InsertXML($updatedXml)\n"}
